# Rename the Pearson/BTEC logo pictures in the headers and footers.
#
# The diff only touches the `name="..."` attribute Word stores on each
# inline picture (the `wp:docPr`/`pic:cNvPr` pair) - the embedded image
# bytes, relationship ids, sizes, etc. all stay the same:
#
#   footer (Pearson logo, PearsonLogo.png)  : image1.png -> image2.png
#   header (BTEC logo, BTec_Logo-Orange)    : image2.jpg -> image1.jpg
#
# InlineShape has no writable "Name" on the Word object model (only
# Title/AlternativeText); the supported way to rename a picture is to
# round-trip it through a floating Shape (which does expose .Name) and
# convert it back to an inline shape afterwards.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlinePicture($range, $index, $newName) {
    $inline = $range.InlineShapes($index)
    $floating = $inline.ConvertToShape()
    $floating.Name = $newName
    [void]$floating.ConvertToInlineShape()
}

# Footers - Pearson logo: image1.png -> image2.png
Rename-InlinePicture $sec.Footers(1).Range 1 "image2.png"
Rename-InlinePicture $sec.Footers(2).Range 1 "image2.png"

# Headers - BTEC logo: image2.jpg -> image1.jpg
Rename-InlinePicture $sec.Headers(1).Range 1 "image1.jpg"
Rename-InlinePicture $sec.Headers(2).Range 1 "image1.jpg"
